$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 97 (shifts old rows 97-101 down to 98-102,
# carrying their B/H/L/N product data down with them).
$ws.Rows("97:97").Insert()

# Copy formatting (styles) from the row below (old row 97, now row 98)
# so the new row matches the product-row look (same style ids).
$ws.Range("A98:N98").Copy()
$ws.Range("A97:N97").PasteSpecial(-4122)

# Re-create the merged cells for the new row (B:G, H:K, L:M), matching every
# other product row's merge pattern.
$ws.Range("B97:G97").Merge()
$ws.Range("H97:K97").Merge()
$ws.Range("L97:M97").Merge()

# Fill in the brand new product row's data.
$ws.Range("A97").Value = 94
$ws.Range("B97").Value = "معجون سيجنال 120 مل "
$ws.Range("H97").Value = "5:0"
$ws.Range("L97").Value = 60
$ws.Range("N97").Value = "1:0"

# Column A is just the running index (row - 3); the insert shifted it along
# with the rest of the row, so restore the correct sequential numbering for
# the rows that moved down.
$ws.Range("A98").Value = 95
$ws.Range("A99").Value = 96
$ws.Range("A100").Value = 97

# Update the totals row (was row 100, now row 101) to include the new row's
# quantity (+60).
$ws.Range("K101").Value = $ws.Range("K101").Value2 + 60

# Match row heights to the report's layout.
$ws.Rows("97:97").RowHeight = 25.5
$ws.Rows("98:98").RowHeight = 25.5
$ws.Rows("99:99").RowHeight = 24.75
$ws.Rows("100:100").RowHeight = 25.5
$ws.Rows("101:101").RowHeight = 25.5
$ws.Rows("102:102").RowHeight = 16.5
